$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.651.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.449.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.99%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("E9").Value = "  -4.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.288.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.95%  "
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.334"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.46%  "
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.902.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "68.515.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("E16").Value = "  -3.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.466.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.577.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0815"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.90%  "
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "436.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.86%  "
$ws.Range("E35").Value = "  -5.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.300"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "37.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.74%  "
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "132.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.03%  "
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.481"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.558"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.50%  "
